# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.724.47'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +3.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.442.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.08%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '184.11'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.632'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.436.71'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +3.08%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('E10').Value = '  +2.47%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.643'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('E12').Value = '  +4.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000278'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.984.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.61'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.437.52'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.86%  '
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '66.790.81'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.39%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.07'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.02'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '492.47'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +8.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '16.57'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +16.70%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.04'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.49%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.24'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.57%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.63'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.09'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.18%  '
$ws.Range('E28').Value = '  +2.37%  '
$ws.Range('E29').Value = '  +5.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.13'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +8.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '11.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '593.52'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.19'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.112'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.09%  '
$ws.Range('E37').Value = '  +4.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.61'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0772'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +4.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.385'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.70%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.21'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.68%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.181.95'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.96%  '
$ws.Range('E43').Value = '  +4.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0429'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.21%  '
$ws.Range('E45').Value = '  +4.36%  '
$ws.Range('E46').Value = '  +22.08%  '
$ws.Range('E47').Value = '  +1.12%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.21'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.78'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +8.22%  '
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '140.29'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.17%  '
